# Updates the worksheet date and every division-problem answer in the
# table to the new set of values, per the commit's regenerated content.

$d = $word.ActiveDocument

$replacements = @(
    @("2025-11-24 Monday", "2025-11-25 Tuesday"),

    @("18÷9=2, 0",  "92÷6=15, 2"),
    @("61÷5=12, 1", "80÷6=13, 2"),
    @("13÷8=1, 5",  "37÷2=18, 1"),
    @("27÷3=9, 0",  "59÷3=19, 2"),
    @("43÷2=21, 1", "43÷8=5, 3"),

    @("15÷9=1, 6",  "97÷9=10, 7"),
    @("37÷9=4, 1",  "15÷6=2, 3"),
    @("11÷6=1, 5",  "20÷5=4, 0"),
    @("25÷6=4, 1",  "66÷5=13, 1"),
    @("89÷4=22, 1", "14÷3=4, 2"),

    @("82÷7=11, 5", "47÷6=7, 5"),
    @("53÷5=10, 3", "36÷2=18, 0"),
    @("19÷7=2, 5",  "67÷6=11, 1"),
    @("32÷2=16, 0", "99÷5=19, 4"),
    @("36÷8=4, 4",  "65÷7=9, 2"),

    @("14÷5=2, 4",  "52÷7=7, 3"),
    @("29÷8=3, 5",  "79÷5=15, 4"),
    @("60÷6=10, 0", "43÷8=5, 3"),
    @("23÷7=3, 2",  "53÷4=13, 1"),
    @("99÷9=11, 0", "36÷8=4, 4"),

    @("61÷2=30, 1", "97÷8=12, 1"),
    @("86÷3=28, 2", "40÷3=13, 1"),
    @("60÷8=7, 4",  "49÷7=7, 0"),
    @("81÷7=11, 4", "64÷5=12, 4"),
    @("10÷2=5, 0",  "95÷4=23, 3")
)

foreach ($pair in $replacements) {
    $oldText = $pair[0]
    $newText = $pair[1]

    $range = $d.Content
    $found = $range.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)

    if (-not $found) {
        Write-Host "WARNING: text not found -> " $oldText
    }
}

$d.Save()
